$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.652.24"
$ws.Cells.Item(2, 5).Value = "  +0.79%  "
$ws.Cells.Item(3, 4).Value = "1.565.91"
$ws.Cells.Item(3, 5).Value = "  -0.48%  "
$ws.Cells.Item(4, 4).Formula = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Formula = "'210.57"
$ws.Cells.Item(5, 5).Value = "  -0.72%  "
$ws.Cells.Item(6, 4).Formula = "'0.487"
$ws.Cells.Item(6, 5).Value = "  -0.89%  "
$ws.Cells.Item(7, 5).Value = "  -0.14%  "
$ws.Cells.Item(8, 4).Formula = "'25.00"
$ws.Cells.Item(8, 5).Value = "  +5.34%  "
$ws.Cells.Item(9, 5).Value = "  -0.76%  "
$ws.Cells.Item(10, 5).Value = "  -0.27%  "
$ws.Cells.Item(11, 4).Formula = "'0.0896"
$ws.Cells.Item(11, 5).Value = "  +0.25%  "
$ws.Cells.Item(12, 4).Value = "1.789.63"
$ws.Cells.Item(12, 5).Value = "  -0.47%  "
$ws.Cells.Item(13, 4).Value = "1.568.13"
$ws.Cells.Item(13, 5).Value = "  -0.38%  "
$ws.Cells.Item(14, 4).Value = "28.670.70"
$ws.Cells.Item(14, 5).Value = "  +0.90%  "
$ws.Cells.Item(15, 5).Value = "  -0.42%  "
$ws.Cells.Item(16, 5).Value = "  -1.59%  "
$ws.Cells.Item(17, 5).Value = "  -0.43%  "
$ws.Cells.Item(18, 4).Formula = "'232.05"
$ws.Cells.Item(18, 5).Value = "  +0.84%  "
$ws.Cells.Item(19, 4).Formula = "'7.38"
$ws.Cells.Item(19, 5).Value = "  -0.32%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0676"
$ws.Cells.Item(20, 5).Value = "  -1.26%  "
$ws.Cells.Item(21, 5).Value = "  -0.10%  "
$ws.Cells.Item(22, 5).Value = "  -1.32%  "
$ws.Cells.Item(23, 4).Formula = "'8.97"
$ws.Cells.Item(23, 5).Value = "  -0.52%  "
$ws.Cells.Item(24, 4).Formula = "'2.09"
$ws.Cells.Item(24, 5).Value = "  +2.00%  "
$ws.Cells.Item(25, 4).Formula = "'150.22"
$ws.Cells.Item(26, 4).Formula = "'14.78"
$ws.Cells.Item(26, 5).Value = "  -0.88%  "
$ws.Cells.Item(27, 5).Value = "  -0.15%  "
$ws.Cells.Item(28, 5).Value = "  -0.07%  "
$ws.Cells.Item(29, 5).Value = "  -2.25%  "
$ws.Cells.Item(30, 4).Formula = "'0.0461"
$ws.Cells.Item(30, 5).Value = "  -4.73%  "
$ws.Cells.Item(31, 5).Value = "  -1.74%  "
$ws.Cells.Item(32, 5).Value = "  -0.59%  "
$ws.Cells.Item(33, 4).Value = "1.391.98"
$ws.Cells.Item(33, 5).Value = "  +0.78%  "
$ws.Cells.Item(34, 5).Value = "  -4.33%  "
$ws.Cells.Item(35, 5).Value = "  -2.86%  "
$ws.Cells.Item(37, 4).Formula = "'2.68"
$ws.Cells.Item(37, 5).Value = "  +1.13%  "
$ws.Cells.Item(38, 5).Value = "  -2.63%  "
$ws.Cells.Item(39, 4).Formula = "'0.0162"
$ws.Cells.Item(39, 5).Value = "  -1.02%  "
$ws.Cells.Item(40, 5).Value = "  +2.39%  "
$ws.Cells.Item(41, 5).Value = "  -0.73%  "
$ws.Cells.Item(42, 4).Formula = "'0.999"
$ws.Cells.Item(42, 5).Value = "  -0.07%  "
$ws.Cells.Item(43, 4).Formula = "'0.774"
$ws.Cells.Item(43, 5).Value = "  -1.89%  "
$ws.Cells.Item(44, 5).Value = "  -3.02%  "
$ws.Cells.Item(45, 4).Formula = "'63.85"
$ws.Cells.Item(45, 5).Value = "  +2.35%  "
$ws.Cells.Item(46, 4).Formula = "'5.24"
$ws.Cells.Item(46, 5).Value = "  -2.29%  "
$ws.Cells.Item(47, 4).Value = "1.701.79"
$ws.Cells.Item(47, 5).Value = "  -0.50%  "
$ws.Cells.Item(48, 4).Formula = "'0.869"
$ws.Cells.Item(48, 5).Value = "  -5.42%  "
$ws.Cells.Item(49, 4).Formula = "'85.26"
$ws.Cells.Item(49, 5).Value = "  +0.03%  "
$ws.Cells.Item(50, 4).Formula = "'43.28"
$ws.Cells.Item(50, 5).Value = "  +6.49%  "
$ws.Cells.Item(51, 4).Value = "0.0₆0100"
$ws.Cells.Item(51, 5).Value = "  -0.10%  "
